$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "39.529.59"
$ws.Range("E2").Value = "  -5.03%  "

# Row 3
$ws.Range("D3").Value = "2.311.32"
$ws.Range("E3").Value = "  -6.02%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "'305.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.20%  "

# Row 6
$ws.Range("D6").Value = "'83.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.30%  "

# Row 7
$ws.Range("E7").Value = "  -3.67%  "

# Row 8
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("E9").Value = "  -4.98%  "

# Row 10
$ws.Range("D10").Value = "'0.0804"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.85%  "

# Row 11
$ws.Range("D11").Value = "'29.47"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -9.32%  "

# Row 12
$ws.Range("E12").Value = "  +0.18%  "

# Row 13
$ws.Range("D13").Value = "2.671.92"
$ws.Range("E13").Value = "  -5.91%  "

# Row 14
$ws.Range("D14").Value = "'6.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.07%  "

# Row 15
$ws.Range("E15").Value = "  -5.60%  "

# Row 16
$ws.Range("D16").Value = "2.315.11"
$ws.Range("E16").Value = "  -6.83%  "

# Row 17
$ws.Range("D17").Value = "'0.746"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.76%  "

# Row 18
$ws.Range("D18").Value = "39.548.54"
$ws.Range("E18").Value = "  -4.76%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0892"
$ws.Range("E19").Value = "  -4.65%  "

# Row 20
$ws.Range("E20").Value = "  -5.19%  "

# Row 21
$ws.Range("D21").Value = "'67.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.12%  "

# Row 22
$ws.Range("E22").Value = "  -5.53%  "

# Row 23
$ws.Range("D23").Value = "'234.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.23%  "

# Row 24
$ws.Range("E24").Value = "  -7.81%  "

# Row 26
$ws.Range("E26").Value = "  -7.36%  "

# Row 27
$ws.Range("D27").Value = "'22.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.79%  "

# Row 28
$ws.Range("E28").Value = "  -5.37%  "

# Row 29
$ws.Range("D29").Value = "'9.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.81%  "

# Row 30
$ws.Range("D30").Value = "'33.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.97%  "

# Row 31
$ws.Range("D31").Value = "'149.94"
$ws.Range("D31").Style = "Normal"

# Row 33
$ws.Range("E33").Value = "  -6.42%  "

# Row 34
$ws.Range("E34").Value = "  -4.85%  "

# Row 35
$ws.Range("D35").Value = "'0.0712"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.97%  "

# Row 36
$ws.Range("E36").Value = "  -2.51%  "

# Row 37
$ws.Range("D37").Value = "'0.0982"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.23%  "

# Row 38
$ws.Range("D38").Value = "'2.70"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.41%  "

# Row 39
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").Value = "'15.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -9.51%  "

# Row 40
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "'1.68"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.23%  "

# Row 41
$ws.Range("E41").Value = "  -6.00%  "

# Row 42
$ws.Range("D42").Value = "'2.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.48%  "

# Row 43
$ws.Range("D43").Value = "1.929.81"
$ws.Range("E43").Value = "  -2.99%  "

# Row 44
$ws.Range("E44").Value = "  -6.72%  "

# Row 45
$ws.Range("D45").Value = "'17.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.04%  "

# Row 46
$ws.Range("E46").Value = "  -1.96%  "

# Row 47
$ws.Range("E47").Value = "  -9.67%  "

# Row 48
$ws.Range("D48").Value = "2.538.82"
$ws.Range("E48").Value = "  -6.55%  "

# Row 49
$ws.Range("D49").Value = "'91.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.17%  "

# Row 50
$ws.Range("D50").Value = "'69.08"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.22%  "

# Row 51
$ws.Range("D51").Value = "'62.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.78%  "
